$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '20.111.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -7.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.430.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -7.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '275.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3744'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3092'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '40.16'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.014'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06605'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.15%  '
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.409'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.26'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.187'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.429.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.16%  '
$ws.Range('E17').Value = '  -8.56%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '75.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -9.13%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.05822'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -11.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.696'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.54'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.336'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '20.113.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.296'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '138.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.591.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '109.56'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.964'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -18.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9148'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.433'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07776'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.11%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.448'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.84%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.54%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05706'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.03%  '
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.000'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.766'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.1926'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.123'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.02029'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.300'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5352'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.549'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5153'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.782'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '109.88'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.67%  '
$ws.Range('E50').Value = '  -6.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.16%  '
